$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at A; existing columns A:E (headers A,B,C,D,F)
# shift right to B:F.
$ws.Columns("A").Insert()

# New column A should look like the other header cells (bold/bordered/
# centered style) -> copy the format from the adjacent header cell.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# New header label for the inserted ID column.
$ws.Range("A1").Value = "ID"

# Row labels (IDs) for each data row, rows 2-25.
$ids = @(
    "Hb 2", "Hb 3", "S 24", "S 28", "Hb 107", "Hb 66", "Hb 69", "Hb 95",
    "Hb 99", "Hb 92", "Hb 40", "Hb 41", "S 11", "Hb 57", "S 21", "S 22",
    "S 3", "S 4", "S 5", "Hb 74", "Hb 79", "Hb 32", "S 15", "S 16"
)

$r = 2
foreach ($id in $ids) {
    $ws.Cells.Item($r, 1).Value = $id
    $r = $r + 1
}
